$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H19").Value = 2147.8
$ws.Range("I19").Value = 2049.8
$ws.Range("J19").Value = 2245.8
$ws.Range("K19").Value = 2049.8
$ws.Range("L19").Value = 2245.8
$ws.Range("M19").Value = -1874.8
$ws.Range("N19").Value = -2595.8
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H53").Value = 269.5
$ws.Range("I53").Value = 366.875
$ws.Range("J53").Value = 74.75
$ws.Range("K53").Value = 366.875
$ws.Range("L53").Value = 74.75
$ws.Range("M53").Value = 270.125
$ws.Range("N53").Value = -1348.75
$ws.Range("H64").Value = 4990.5454
$ws.Range("I64").Value = 4990.5454
$ws.Range("K64").Value = 4990.5454
$ws.Range("M64").Value = -4742.5454
$ws.Range("H67").Value = 4990.5454
$ws.Range("I67").Value = 4990.5454
$ws.Range("K67").Value = 4990.5454
$ws.Range("M67").Value = -4132.5454
$ws.Range("H106").Value = 2109.3635
$ws.Range("I106").Value = 2210.3
$ws.Range("K106").Value = 2210.3
$ws.Range("M106").Value = -1579.3
$ws.Range("H135").Value = 950
$ws.Range("I135").Value = 950
$ws.Range("K135").Value = 8550
$ws.Range("M135").Value = -6015
$ws.Range("H137").Value = 3669
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3669
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 11007
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -16107
$ws.Range("H138").Value = 2727
$ws.Range("I138").Value = 1909.2858
$ws.Range("J138").Value = 2999.5715
$ws.Range("K138").Value = 5727.857400000001
$ws.Range("L138").Value = 8998.7145
$ws.Range("M138").Value = -587.8574000000008
$ws.Range("N138").Value = -19278.7145
$ws.Range("H140").Value = 20000
$ws.Range("J140").Value = 20000
$ws.Range("L140").Value = 20000
$ws.Range("N140").Value = -30360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 26000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 26000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 26000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -26518
$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 5000
$ws.Range("K37").Value = 5000
$ws.Range("M37").Value = -4727
$ws.Range("H55").Value = 31666.666
$ws.Range("J55").Value = 31666.666
$ws.Range("L55").Value = 31666.666
$ws.Range("N55").Value = -32296.666
$ws.Range("H74").Value = 2324.6667
$ws.Range("I74").Value = 1955.75
$ws.Range("J74").Value = 3062.5
$ws.Range("K74").Value = 1955.75
$ws.Range("L74").Value = 3062.5
$ws.Range("M74").Value = -1081.75
$ws.Range("N74").Value = -4810.5
$ws.Range("H77").Value = 2324.6667
$ws.Range("I77").Value = 1955.75
$ws.Range("J77").Value = 3062.5
$ws.Range("K77").Value = 9778.75
$ws.Range("L77").Value = 15312.5
$ws.Range("M77").Value = -5410.75
$ws.Range("N77").Value = -24048.5
$ws.Range("H88").Value = 4147.3335
$ws.Range("I88").Value = 2280
$ws.Range("J88").Value = 4380.75
$ws.Range("K88").Value = 2280
$ws.Range("L88").Value = 4380.75
$ws.Range("M88").Value = -1874
$ws.Range("N88").Value = -5192.75
$ws.Range("H91").Value = 4147.3335
$ws.Range("I91").Value = 2280
$ws.Range("J91").Value = 4380.75
$ws.Range("K91").Value = 2280
$ws.Range("L91").Value = 4380.75
$ws.Range("M91").Value = -876
$ws.Range("N91").Value = -7188.75
$ws.Range("H97").Value = 972.2727
$ws.Range("I97").Value = 909.5
$ws.Range("K97").Value = 909.5
$ws.Range("M97").Value = -413.5
$ws.Range("H110").Value = 1579.4
$ws.Range("I110").Value = 1579.4
$ws.Range("K110").Value = 1579.4
$ws.Range("M110").Value = 465.5999999999999
$ws.Range("H132").Value = 1009.8929
$ws.Range("I132").Value = 1019.5185
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 3058.5555
$ws.Range("L132").Value = 2250
$ws.Range("M132").Value = -528.5554999999999
$ws.Range("N132").Value = -7310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6760.3335
$ws.Range("I22").Value = 114.6
$ws.Range("K22").Value = 114.6
$ws.Range("M22").Value = 58.40000000000001
$ws.Range("H64").Value = 1187.0834
$ws.Range("I64").Value = 1435.6
$ws.Range("J64").Value = 1009.5714
$ws.Range("K64").Value = 1435.6
$ws.Range("L64").Value = 1009.5714
$ws.Range("M64").Value = -1210.6
$ws.Range("N64").Value = -1459.5714
$ws.Range("H67").Value = 1187.0834
$ws.Range("I67").Value = 1435.6
$ws.Range("J67").Value = 1009.5714
$ws.Range("K67").Value = 1435.6
$ws.Range("L67").Value = 1009.5714
$ws.Range("M67").Value = -655.5999999999999
$ws.Range("N67").Value = -2569.5714
$ws.Range("H107").Value = 1144.4286
$ws.Range("I107").Value = 1001.8333
$ws.Range("K107").Value = 1001.8333
$ws.Range("M107").Value = 918.1667
$ws.Range("H134").Value = 3292.2222
$ws.Range("I134").Value = 3309.4119
$ws.Range("K134").Value = 9928.235700000001
$ws.Range("M134").Value = -7393.235700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3468.25
$ws.Range("I31").Value = 2511
$ws.Range("J31").Value = 4425.5
$ws.Range("K31").Value = 2511
$ws.Range("L31").Value = 4425.5
$ws.Range("M31").Value = -2216
$ws.Range("N31").Value = -5015.5
$ws.Range("H34").Value = 3468.25
$ws.Range("I34").Value = 2511
$ws.Range("J34").Value = 4425.5
$ws.Range("K34").Value = 2511
$ws.Range("L34").Value = 4425.5
$ws.Range("M34").Value = -2309
$ws.Range("N34").Value = -4829.5
$ws.Range("H64").Value = 28000
$ws.Range("J64").Value = 28000
$ws.Range("L64").Value = 28000
$ws.Range("N64").Value = -28496
$ws.Range("H67").Value = 28000
$ws.Range("J67").Value = 28000
$ws.Range("L67").Value = 28000
$ws.Range("N67").Value = -29716
$ws.Range("H134").Value = 2642.5293
$ws.Range("I134").Value = 2338.2727
$ws.Range("J134").Value = 3200.3333
$ws.Range("K134").Value = 7014.8181
$ws.Range("L134").Value = 9600.999899999999
$ws.Range("M134").Value = -4479.8181
$ws.Range("N134").Value = -14670.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 37.166668
$ws.Range("I12").Value = 26.5
$ws.Range("J12").Value = 47.833332
$ws.Range("K12").Value = 79.5
$ws.Range("L12").Value = 143.499996
$ws.Range("M12").Value = 93.5
$ws.Range("N12").Value = -489.499996
$ws.Range("H23").Value = 218.25
$ws.Range("I23").Value = 19
$ws.Range("K23").Value = 57
$ws.Range("M23").Value = 178
$ws.Range("H36").Value = 325
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H60").Value = 513.8570999999999
$ws.Range("I60").Value = 319.4
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 958.1999999999999
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -707.1999999999999
$ws.Range("N60").Value = -3502
$ws.Range("H92").Value = 126
$ws.Range("I92").Value = 126
$ws.Range("K92").Value = 378
$ws.Range("M92").Value = 870
$ws.Range("H107").Value = 685.2222
$ws.Range("J107").Value = 657.3077
$ws.Range("L107").Value = 1971.9231
$ws.Range("N107").Value = -5811.9231
$ws.Range("H121").Value = 824.9375
$ws.Range("I121").Value = 703.6667
$ws.Range("J121").Value = 980.8570999999999
$ws.Range("K121").Value = 2111.0001
$ws.Range("L121").Value = 2942.5713
$ws.Range("M121").Value = -801.0001000000002
$ws.Range("N121").Value = -5562.5713
$ws.Range("H122").Value = 455.75
$ws.Range("J122").Value = 292
$ws.Range("L122").Value = 2628
$ws.Range("N122").Value = -7528
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3594
$ws.Range("I80").Value = 3200
$ws.Range("J80").Value = 3988
$ws.Range("K80").Value = 3200
$ws.Range("L80").Value = 3988
$ws.Range("M80").Value = -2202
$ws.Range("N80").Value = -5984
$ws.Range("H83").Value = 3594
$ws.Range("I83").Value = 3200
$ws.Range("J83").Value = 3988
$ws.Range("K83").Value = 16000
$ws.Range("L83").Value = 19940
$ws.Range("M83").Value = -11008
$ws.Range("N83").Value = -29924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1300
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1600
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1600
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2190
$ws.Range("H27").Value = 1300
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1600
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1600
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1814
$ws.Range("H46").Value = 3574.6365
$ws.Range("I46").Value = 2626.889
$ws.Range("J46").Value = 4230.769
$ws.Range("K46").Value = 2626.889
$ws.Range("L46").Value = 4230.769
$ws.Range("M46").Value = -2438.889
$ws.Range("N46").Value = -4606.769
$ws.Range("H136").Value = 2493.4285
$ws.Range("I136").Value = 2584.1667
$ws.Range("J136").Value = 1949
$ws.Range("K136").Value = 7752.500100000001
$ws.Range("L136").Value = 5847
$ws.Range("M136").Value = -5202.500100000001
$ws.Range("N136").Value = -10947

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3571
$ws.Range("I81").Value = 3571
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7142
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6081
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 3571
$ws.Range("I84").Value = 3571
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 35710
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -30406
$ws.Range("N84").ClearContents()
Write-Host "Edits applied successfully"
